$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOGT1")

# New row 13: continuation of the time log, and a new Phase/Task comment.
$ws.Range("A13").Value = 41931
$ws.Range("B13").Value = 0.70972222222222225
$ws.Range("C13").Value = 0.74305555555555547
$ws.Range("D13").Value = 0
$ws.Range("E13").Formula = "=((HOUR(C13)-HOUR(B13))*60)+(MINUTE(C13)-MINUTE(B13))-D13"
$ws.Range("F13").Value = 43
$ws.Range("H13").Value = "Presentar al equipo herramienta Rails"

$ws.Range("A13:H13").RowHeight = 52

$ws.Range("H14").Select()
